$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "28-10-2025"
$ws.Range("B44").Value = "The price of gold in India today is ₹12,246 per gram for 24 karat gold, ₹11,225 per gram for 22 karat gold and ₹9,184 per gram for 18 karat gold (also called 999 gold)."
